$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Diseases (patient-stated)
$ws.Range("B2").Value = 5.9
$ws.Range("C2").Value = 2.9

# Row 3: Injuries & adverse effects
$ws.Range("B3").Value = 14.9
$ws.Range("C3").Value = 8.300000000000001
$ws.Range("D3").Value = 6.7

# Row 4: Other
$ws.Range("B4").Value = 6.8
$ws.Range("C4").Value = 5.2

# Row 5: Symptom - Circulatory
$ws.Range("B5").Value = 10.1
$ws.Range("C5").Value = 6.4

# Row 6: Symptom - Digestive
$ws.Range("B6").Value = 12.9
$ws.Range("C6").Value = 7.5
$ws.Range("D6").Value = 6.3

# Row 7: Symptom - General
$ws.Range("B7").Value = 4.8
$ws.Range("C7").Value = 3.6
$ws.Range("D7").Value = 5.9

# Row 8: Symptom - Genitourinary -> Symptom - Musculoskeletal
$ws.Range("A8").Value = "Symptom – Musculoskeletal"
$ws.Range("B8").Value = 2.5
$ws.Range("C8").Value = 1.1
$ws.Range("D8").Value = 2.4

# Row 9: Symptom - Nervous
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = 11.5
$ws.Range("D9").Value = 10.4

# Row 10: Symptom - Respiratory
$ws.Range("B10").Value = 26.3
$ws.Range("C10").Value = 51.3
$ws.Range("D10").Value = 44.8

# Row 11: Symptom - Skin/Hair/Nails
$ws.Range("B11").Value = 2.7
$ws.Range("C11").Value = 1.3
$ws.Range("D11").Value = 2.8

# Row 12: Uncodable/Unknown
$ws.Range("B12").Value = 2.2
$ws.Range("C12").Value = 0.9
